$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 corresponds to theta_se; Row 6 corresponds to lambda_se.
# Previously every cell in these rows held the shared string "(nan)".
# They now hold distinct standard-error values (output of a pickled
# multiple-imputation run), one unique string per cell. Values are
# written column-by-column (theta_se then lambda_se for each column)
# to match the order new shared-string entries were created.

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L")

$row4Values = @("(0.22)", "(0.63)", "(0.3)", "(0.18)", "(2.41)", "(1.47)", "(0.75)", "(0.25)", "(1.03)", "(2.13)", "(2.36)")
$row6Values = @("(0.12)", "(0.4)", "(0.2)", "(0.19)", "(0.35)", "(1.09)", "(0.71)", "(0.11)", "(0.92)", "(1.29)", "(2.03)")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $ws.Range("$col`4").Value = $row4Values[$i]
    $ws.Range("$col`6").Value = $row6Values[$i]
}
